$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Extend the table with a new row 54, continuing the same pattern as row 53.
# Copy A53's format (date number format + border) onto A54 before writing the value.
$ws.Range("A53").Copy($ws.Range("A54"))

$ws.Range("A54").Value = 45986
$ws.Range("B54").Value = 2025
$ws.Range("C54").Value = -2.451276118722334
$ws.Range("D54").Value = 2026
$ws.Range("E54").Value = 1.795477855501626

$wb.Save()
